$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1971.9166
$ws.Range("I11").Value = 1971.9166
$ws.Range("K11").Value = 1971.9166
$ws.Range("M11").Value = -1831.9166
$ws.Range("H32").Value = 4332.3335
$ws.Range("I32").Value = 4332.3335
$ws.Range("K32").Value = 4332.3335
$ws.Range("M32").Value = -4006.3335
$ws.Range("H40").Value = 14288700
$ws.Range("I40").Value = 3149.2856
$ws.Range("J40").Value = 21431476
$ws.Range("K40").Value = 3149.2856
$ws.Range("L40").Value = 21431476
$ws.Range("M40").Value = -2974.2856
$ws.Range("N40").Value = -21431826
$ws.Range("H43").Value = 3660.8333
$ws.Range("J43").Value = 6681.25
$ws.Range("L43").Value = 6681.25
$ws.Range("N43").Value = -6819.25
$ws.Range("H69").Value = 9573.571
$ws.Range("I69").Value = 9573.571
$ws.Range("K69").Value = 28720.713
$ws.Range("M69").Value = -27846.713
$ws.Range("H72").Value = 9573.571
$ws.Range("I72").Value = 9573.571
$ws.Range("K72").Value = 86162.139
$ws.Range("M72").Value = -81794.139
$ws.Range("H86").Value = 3397.0588
$ws.Range("I86").Value = 744.1667
$ws.Range("J86").Value = 9764
$ws.Range("K86").Value = 744.1667
$ws.Range("L86").Value = 9764
$ws.Range("M86").Value = 378.8333
$ws.Range("N86").Value = -12010
$ws.Range("H89").Value = 3397.0588
$ws.Range("I89").Value = 744.1667
$ws.Range("J89").Value = 9764
$ws.Range("K89").Value = 3720.8335
$ws.Range("L89").Value = 48820
$ws.Range("M89").Value = 1895.1665
$ws.Range("N89").Value = -60052
$ws.Range("H98").Value = 1671067.5
$ws.Range("I98").Value = 5281
$ws.Range("K98").Value = 5281
$ws.Range("M98").Value = -3783
$ws.Range("H122").Value = 1671067.5
$ws.Range("I122").Value = 5281
$ws.Range("K122").Value = 15843
$ws.Range("M122").Value = -13393
$ws.Range("H132").Value = 2996.262
$ws.Range("I132").Value = 2972.147
$ws.Range("K132").Value = 8916.440999999999
$ws.Range("M132").Value = -6386.440999999999
$ws.Range("H137").Value = 4047324
$ws.Range("I137").Value = 9192880
$ws.Range("J137").Value = 4387.5713
$ws.Range("K137").Value = 27578640
$ws.Range("L137").Value = 13162.7139
$ws.Range("M137").Value = -27576090
$ws.Range("N137").Value = -18262.7139
$ws.Range("H141").Value = 2298.2354
$ws.Range("I141").Value = 1965.6666
$ws.Range("K141").Value = 5896.9998
$ws.Range("M141").Value = -716.9997999999996

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4552.381
$ws.Range("I61").Value = 3203.8667
$ws.Range("K61").Value = 3203.8667
$ws.Range("M61").Value = -2991.8667
$ws.Range("H96").Value = 7342
$ws.Range("J96").Value = 7342
$ws.Range("L96").Value = 7342
$ws.Range("N96").Value = -12834
$ws.Range("H122").Value = 4094.5
$ws.Range("I122").Value = 4012
$ws.Range("J122").Value = 4507
$ws.Range("K122").Value = 12036
$ws.Range("L122").Value = 13521
$ws.Range("M122").Value = -9586
$ws.Range("N122").Value = -18421
$ws.Range("H136").Value = 4552.381
$ws.Range("I136").Value = 3203.8667
$ws.Range("K136").Value = 9611.6001
$ws.Range("M136").Value = -7061.6001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1693.909
$ws.Range("I80").Value = 2439.6
$ws.Range("J80").Value = 1072.5
$ws.Range("K80").Value = 2439.6
$ws.Range("L80").Value = 1072.5
$ws.Range("M80").Value = -1441.6
$ws.Range("N80").Value = -3068.5
$ws.Range("H83").Value = 1693.909
$ws.Range("I83").Value = 2439.6
$ws.Range("J83").Value = 1072.5
$ws.Range("K83").Value = 12198
$ws.Range("L83").Value = 5362.5
$ws.Range("M83").Value = -7206
$ws.Range("N83").Value = -15346.5
$ws.Range("H107").Value = 2261.85
$ws.Range("I107").Value = 2361.1765
$ws.Range("J107").Value = 1699
$ws.Range("K107").Value = 2361.1765
$ws.Range("L107").Value = 1699
$ws.Range("M107").Value = -441.1765
$ws.Range("N107").Value = -5539

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 557997.75
$ws.Range("I31").Value = 835451.5600000001
$ws.Range("K31").Value = 835451.5600000001
$ws.Range("M31").Value = -835156.5600000001
$ws.Range("H34").Value = 557997.75
$ws.Range("I34").Value = 835451.5600000001
$ws.Range("K34").Value = 835451.5600000001
$ws.Range("M34").Value = -835249.5600000001
$ws.Range("H132").Value = 1971.8334
$ws.Range("I132").Value = 1950.7
$ws.Range("K132").Value = 5852.1
$ws.Range("M132").Value = -3322.1

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null
$ws.Range("H97").Value = 999.6667
$ws.Range("J97").Value = 999.6667
$ws.Range("L97").Value = 2999.0001
$ws.Range("N97").Value = -3991.0001
$ws.Range("H122").Value = 907.2857
$ws.Range("J122").Value = 939.6
$ws.Range("L122").Value = 8456.4
$ws.Range("N122").Value = -13356.4
$ws.Range("H131").Value = 14287062
$ws.Range("I131").Value = 62500724
$ws.Range("J131").Value = 1532.4073
$ws.Range("K131").Value = 187502172
$ws.Range("L131").Value = 4597.2219
$ws.Range("M131").Value = -187497132
$ws.Range("N131").Value = -14677.2219

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4112.0303
$ws.Range("I80").Value = 3652.55
$ws.Range("J80").Value = 4818.923
$ws.Range("K80").Value = 3652.55
$ws.Range("L80").Value = 4818.923
$ws.Range("M80").Value = -2654.55
$ws.Range("N80").Value = -6814.923
$ws.Range("H83").Value = 4112.0303
$ws.Range("I83").Value = 3652.55
$ws.Range("J83").Value = 4818.923
$ws.Range("K83").Value = 18262.75
$ws.Range("L83").Value = 24094.615
$ws.Range("M83").Value = -13270.75
$ws.Range("N83").Value = -34078.615
$ws.Range("H122").Value = 5323.5
$ws.Range("I122").Value = 4989.2144
$ws.Range("K122").Value = 14967.6432
$ws.Range("M122").Value = -12517.6432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1654.8718
$ws.Range("I16").Value = 1638.3793
$ws.Range("J16").Value = 1702.7
$ws.Range("K16").Value = 1638.3793
$ws.Range("L16").Value = 1702.7
$ws.Range("M16").Value = -1468.3793
$ws.Range("N16").Value = -2042.7
$ws.Range("H61").Value = 2802.111
$ws.Range("I61").Value = 2802.111
$ws.Range("K61").Value = 2802.111
$ws.Range("M61").Value = -2600.111
$ws.Range("H68").Value = 8637.526
$ws.Range("I68").Value = 9006.277
$ws.Range("K68").Value = 9006.277
$ws.Range("M68").Value = -8257.277
$ws.Range("H71").Value = 8637.526
$ws.Range("I71").Value = 9006.277
$ws.Range("K71").Value = 45031.385
$ws.Range("M71").Value = -41287.385
$ws.Range("H82").Value = 5324.75
$ws.Range("I82").Value = 2300
$ws.Range("J82").Value = 6333
$ws.Range("K82").Value = 2300
$ws.Range("L82").Value = 6333
$ws.Range("M82").Value = -1939
$ws.Range("N82").Value = -7055
$ws.Range("H85").Value = 5324.75
$ws.Range("I85").Value = 2300
$ws.Range("J85").Value = 6333
$ws.Range("K85").Value = 2300
$ws.Range("L85").Value = 6333
$ws.Range("M85").Value = -1052
$ws.Range("N85").Value = -8829
$ws.Range("H104").Value = 22446.334
$ws.Range("J104").Value = 22446.334
$ws.Range("L104").Value = 22446.334
$ws.Range("N104").Value = -29434.334
$ws.Range("H113").Value = 2802.111
$ws.Range("I113").Value = 2802.111
$ws.Range("K113").Value = 2802.111
$ws.Range("M113").Value = -632.1109999999999
$ws.Range("H132").Value = 3277.6296
$ws.Range("I132").Value = 3163.0417
$ws.Range("K132").Value = 9489.125100000001
$ws.Range("M132").Value = -6959.125100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8239.25
$ws.Range("I62").Value = 3430
$ws.Range("K62").Value = 3430
$ws.Range("M62").Value = -2806
$ws.Range("H65").Value = 8239.25
$ws.Range("I65").Value = 3430
$ws.Range("K65").Value = 17150
$ws.Range("M65").Value = -14030
$ws.Range("H107").Value = 560.6667
$ws.Range("I107").Value = 560.6667
$ws.Range("K107").Value = 1682.0001
$ws.Range("M107").Value = 237.9999
$ws.Range("H126").Value = 4811.375
$ws.Range("I126").Value = 4665.3335
$ws.Range("J126").Value = 4899
$ws.Range("K126").Value = 13996.0005
$ws.Range("L126").Value = 14697
$ws.Range("M126").Value = -11526.0005
$ws.Range("N126").Value = -19637
$ws.Range("H132").Value = 1499.5
$ws.Range("I132").Value = 1499
$ws.Range("K132").Value = 4497
$ws.Range("M132").Value = -1967
$ws.Range("H136").Value = 669022.9
$ws.Range("I136").Value = 770745.7
$ws.Range("K136").Value = 2312237.1
$ws.Range("M136").Value = -2309687.1
